$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New log entries added to the "Description" table (new shared strings)
$ws.Range("E5").Value = "Chapter 7 contains new information. I have not worked with dictionaries before. I finished the lessons and started working on practice question 1. But have not finished it. I suspect the lessons will take longer from here on out."
$ws.Range("E6").Value = "Finished practice programs ch. 7. The first one was the most complex, the following two were a lot simpler. I also took longer on the first program because I challenged myself to write more efficient code. Instead of separating the code and looping over a dictionary multiple times I did it in one loop. It may however be better to separate into functions for readability, I will focus on this with the next practice programs."
$ws.Range("E7").Value = "Finished Chapter 8. This was the final chapter of the introductory chapters. It was focused on string manipulation. Next the focus will shift towards automation tasks, which is what I am doing this course for, so I look forward to it."

# Row 5: 2025-09-24, 14:30 - 17:30 (serial date 45924)
$ws.Range("A5").Value = 45924
$ws.Range("B5").Value = 0.60416666666666663
$ws.Range("C5").Value = 0.72916666666666663

# Row 6: 2025-10-01, 10:00 - 12:00 (serial date 45931)
$ws.Range("A6").Value = 45931
$ws.Range("B6").Value = 0.41666666666666669
$ws.Range("C6").Value = 0.5

# Row 7: 2025-10-01, 12:00 - 16:00 (serial date 45931)
$ws.Range("A7").Value = 45931
$ws.Range("B7").Value = 0.5
$ws.Range("C7").Value = 0.66666666666666663

# Row heights to match the diff
$ws.Rows.Item(5).RowHeight = 43.2
$ws.Rows.Item(6).RowHeight = 72
$ws.Rows.Item(7).RowHeight = 43.2

# Update view: scroll so row 3 is the top-left visible row, then select E8
$excel.Goto($ws.Range("A3"), $true)
$ws.Range("E8").Select()

$wb.Save()
